# Update gh-pages output (generated at a56beed)
# Sheet 1 = 展览 (Exhibitions), Sheet 2 = 演出 (Performances),
# Sheet 3 = 本地生活 (Local life), Sheet 4 = 全部类型 (All types)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (index 1): refresh "想去人数" (interest count) column F
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item(1)
$expoUpdates = @{
    2  = 1140
    4  = 11178
    5  = 1434
    6  = 410
    7  = 682
    8  = 2159
    9  = 635
    10 = 838
    11 = 435
    12 = 309
    13 = 355
    14 = 335
    15 = 1079
    16 = 468
    17 = 852
    18 = 305
    19 = 499
    20 = 780
    21 = 834
    22 = 56
    23 = 84
    24 = 216
    25 = 524
    26 = 55
    27 = 25
    28 = 249
}
foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Range("F$row").Value = $expoUpdates[$row]
}

# ---------------------------------------------------------------------
# Sheet "演出" (index 2): all events expired -> remove data rows,
# keep only the header row (dimension collapses to A1:J1)
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item(2)
$wsShow.Rows("2:9").Delete()

# ---------------------------------------------------------------------
# Sheet "全部类型" (index 4): refresh "想去人数" (interest count) column F
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item(4)
$allUpdates = @{
    4  = 1140
    5  = 675
    7  = 29
    8  = 11178
    9  = 1434
    10 = 66
    11 = 411
    12 = 682
    13 = 2159
    14 = 635
    15 = 838
    16 = 22
    17 = 435
    18 = 309
    19 = 355
    20 = 335
    21 = 1079
    22 = 468
    23 = 734
    24 = 852
    25 = 305
    26 = 499
    27 = 780
    28 = 834
    29 = 56
    30 = 13
    31 = 84
    32 = 216
    33 = 524
    34 = 55
    35 = 25
    36 = 249
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
